$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "thinBasic_Test_77110_64"

# Update the header label text ("Few random number" -> "Few random numbers")
$ws.Range("A1").Value = "Few random numbers"

# New set of random numbers in column B (B1:B20) -- drives the SIN() formulas in D
$ws.Range("B1").Value = 845
$ws.Range("B2").Value = 761
$ws.Range("B3").Value = 137
$ws.Range("B4").Value = 1486
$ws.Range("B5").Value = 750
$ws.Range("B6").Value = 242
$ws.Range("B7").Value = 1899
$ws.Range("B8").Value = 1220
$ws.Range("B9").Value = 822
$ws.Range("B10").Value = 582
$ws.Range("B11").Value = 273
$ws.Range("B12").Value = 119
$ws.Range("B13").Value = 1481
$ws.Range("B14").Value = 1347
$ws.Range("B15").Value = 1814
$ws.Range("B16").Value = 457
$ws.Range("B17").Value = 494
$ws.Range("B18").Value = 453
$ws.Range("B19").Value = 1179
$ws.Range("B20").Value = 379

# Center-align the SIN() results column and widen it so the values are readable
$ws.Range("D1:D20").HorizontalAlignment = -4108
$ws.Columns.Item(4).ColumnWidth = 24.8
